# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# Adds a new "ODI Batting Extra" worksheet (after "ODI Batting") with the
# scraped extra batting fields, and removes the stray empty B12 cell on
# the "ODI Batting" sheet.

$wb = $excel.ActiveWorkbook
$wsBatting = $wb.Worksheets.Item("ODI Batting")

# --- 1. Drop the vestigial empty inline-string cell at ODI Batting!B12 ---
$wsBatting.Range("B12").ClearContents()

# --- 2. Create the new sheet, positioned right after "ODI Batting" ---
$newSheet = $wb.Worksheets.Add($null, $wsBatting)
$newSheet.Name = "ODI Batting Extra"

# Helper: write a value as TEXT (inline/shared string), regardless of
# whether it looks numeric, without leaving a stray number-format style
# behind on the cell.
function Set-TextCell($cell, [string]$value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Helper: write a value as a genuine NUMBER.
function Set-NumberCell($cell, $value) {
    $cell.Value = $value
}

# --- Header row (bold + bordered, matching the other sheets) ---
$wsBatting.Range("A1:F1").Copy($newSheet.Range("A1"))
$newSheet.Cells.Item(1,1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1,2).Value = "BATTING_POSITION"
$newSheet.Cells.Item(1,3).Value = "NUM_4"
$newSheet.Cells.Item(1,4).Value = "NUM_6"
$newSheet.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# --- Data rows 2-19 ---
# MATCH_CODE (col A) is identical to ODI Batting!D2:D19 -- copy it across
# so it keeps text formatting cleanly.
$wsBatting.Range("D2:D19").Copy($newSheet.Range("A2"))

$rows = @(
    @{ Row=2;  B=3;     C='2';  D='0'; E='20.45%'; F='NO'  },
    @{ Row=3;  B=3;     C='7';  D='0'; E='26.18%'; F='NO'  },
    @{ Row=4;  B=3;     C='17'; D='0'; E='39.62%'; F='YES' },
    @{ Row=5;  B=4;     C='0';  D='0'; E='2.11%';  F='NO'  },
    @{ Row=6;  B=3;     C='0';  D='0'; E='2.83%';  F='NO'  },
    @{ Row=7;  B=$null; C='';  D='';  E='';        F='NO'  },
    @{ Row=8;  B=$null; C='';  D='';  E='';        F='NO'  },
    @{ Row=9;  B=2;     C='0';  D='0'; E='6.10%';  F='NO'  },
    @{ Row=10; B=2;     C='5';  D='0'; E='8.68%';  F='NO'  },
    @{ Row=11; B=2;     C='3';  D='0'; E='7.77%';  F='NO'  },
    @{ Row=12; B=$null; C='';  D='';  E='';        F='NO'  },
    @{ Row=13; B=2;     C='6';  D='0'; E='36.54%'; F='NO'  },
    @{ Row=14; B=$null; C='';  D='';  E='';        F='NO'  },
    @{ Row=15; B=2;     C='13'; D='1'; E='38.70%'; F='YES' },
    @{ Row=16; B=$null; C='';  D='';  E='';        F='NO'  },
    @{ Row=17; B=$null; C='';  D='';  E='';        F='NO'  },
    @{ Row=18; B=2;     C='1';  D='0'; E='6.48%';  F='NO'  },
    @{ Row=19; B=2;     C='12'; D='8'; E='46.78%'; F='NO'  }
)

foreach ($r in $rows) {
    $rowIdx = $r.Row

    if ($null -eq $r.B) {
        Set-TextCell $newSheet.Cells.Item($rowIdx,2) ""
    } else {
        Set-NumberCell $newSheet.Cells.Item($rowIdx,2) $r.B
    }

    Set-TextCell $newSheet.Cells.Item($rowIdx,3) $r.C
    Set-TextCell $newSheet.Cells.Item($rowIdx,4) $r.D
    Set-TextCell $newSheet.Cells.Item($rowIdx,5) $r.E
    Set-TextCell $newSheet.Cells.Item($rowIdx,6) $r.F
}

$wsBatting.Select()
